# Updates cryptos list (GitHub Actions style refresh): new Price / Volume(1h)
# values for rows 2-51, plus a few rank swaps (Cardano<->Toncoin,
# Dai<->Fetch.AI, PancakeSwap<->NEARProtocol) that moved Coin name + Link
# between adjacent rows.
#
# Note: column D holds plain decimal-looking price text (e.g. "0.997",
# "35.10", "1.00") stored as literal text in the workbook, not numbers.
# Assigning such a string straight to .Value would make Excel silently
# coerce it to a real number (dropping significant trailing zeros, e.g.
# "35.10" -> 35.1). A leading apostrophe forces literal text entry for
# those values, exactly like a user typing '35.10 into the cell; values
# that aren't clean numbers (thousand-separated prices like
# "67.885.61") don't need it since Excel can't coerce them anyway.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.885.61"
$ws.Range("E2").Value = "  +0.09%  "

$ws.Range("D3").Value = "3.787.33"
$ws.Range("E3").Value = "  -0.82%  "

$ws.Range("D4").Value = "'0.997"
$ws.Range("E4").Value = "  -0.25%  "

$ws.Range("D5").Value = "'602.38"
$ws.Range("E5").Value = "  -0.53%  "

$ws.Range("D6").Value = "'163.05"
$ws.Range("E6").Value = "  -2.51%  "

$ws.Range("D7").Value = "3.785.48"
$ws.Range("E7").Value = "  -0.89%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").Value = "'0.514"
$ws.Range("E9").Value = "  -1.31%  "

$ws.Range("E10").Value = "  -2.72%  "

$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").Value = "'6.84"
$ws.Range("E11").Value = "  +8.36%  "

$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").Value = "'0.446"
$ws.Range("E12").Value = "  -1.33%  "

$ws.Range("E13").Value = "  -3.77%  "

$ws.Range("D14").Value = "'35.10"
$ws.Range("E14").Value = "  -2.74%  "

$ws.Range("D15").Value = "4.421.63"
$ws.Range("E15").Value = "  -0.79%  "

$ws.Range("D16").Value = "3.782.64"
$ws.Range("E16").Value = "  -0.42%  "

$ws.Range("D17").Value = "67.818.04"
$ws.Range("E17").Value = "  -0.04%  "

$ws.Range("D18").Value = "'18.17"
$ws.Range("E18").Value = "  -1.82%  "

$ws.Range("E19").Value = "  +1.96%  "

$ws.Range("D20").Value = "'7.01"
$ws.Range("E20").Value = "  -1.36%  "

$ws.Range("D21").Value = "'457.59"
$ws.Range("E21").Value = "  -1.19%  "

$ws.Range("D22").Value = "'9.44"
$ws.Range("E22").Value = "  -4.92%  "

$ws.Range("D23").Value = "'0.690"
$ws.Range("E23").Value = "  -1.72%  "

$ws.Range("D24").Value = "'83.05"
$ws.Range("E24").Value = "  -0.49%  "

$ws.Range("E25").Value = "  -5.65%  "

$ws.Range("E26").Value = "  -2.14%  "

$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D27").Value = "'2.07"
$ws.Range("E27").Value = "  -2.02%  "

$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  -0.18%  "

$ws.Range("D29").Value = "'9.88"
$ws.Range("E29").Value = "  -2.00%  "

$ws.Range("D30").Value = "3.936.82"
$ws.Range("E30").Value = "  -0.80%  "

$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'7.21"
$ws.Range("E31").Value = "  -3.01%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'2.59"
$ws.Range("E32").Value = "  -7.67%  "

$ws.Range("E33").Value = "  -2.52%  "

$ws.Range("D34").Value = "'28.91"
$ws.Range("E34").Value = "  -2.59%  "

$ws.Range("E35").Value = "  -0.23%  "

$ws.Range("D36").Value = "'8.92"
$ws.Range("E36").Value = "  -1.99%  "

$ws.Range("E37").Value = "  -1.15%  "

$ws.Range("E38").Value = "  +4.61%  "

$ws.Range("E39").Value = "  -0.44%  "

$ws.Range("E40").Value = "  -2.35%  "

$ws.Range("E41").Value = "  -6.23%  "

$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.09%  "

$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("D44").Value = "'43.71"
$ws.Range("E44").Value = "  +1.26%  "

$ws.Range("E45").Value = "  -2.16%  "

$ws.Range("D46").Value = "'151.86"
$ws.Range("E46").Value = "  +2.11%  "

$ws.Range("E47").Value = "  -2.58%  "

$ws.Range("D48").Value = "'8.27"
$ws.Range("E48").Value = "  -0.98%  "

$ws.Range("E49").Value = "  -2.31%  "

$ws.Range("E50").Value = "  -1.27%  "

$ws.Range("D51").Value = "'26.29"
$ws.Range("E51").Value = "  -7.84%  "
